# Update "想去人数" (want-to-go count) figures for several camp/expo
# events across the "展览", "演出" and "全部类型" sheets, reflecting a
# refreshed data pull (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# 展览 sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 811
$wsExhibit.Range("F4").Value = 1118
$wsExhibit.Range("F11").Value = 512
$wsExhibit.Range("F14").Value = 12753
$wsExhibit.Range("F16").Value = 5246

# 演出 sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 100

# 全部类型 sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 811
$wsAll.Range("F5").Value = 1118
$wsAll.Range("F12").Value = 512
$wsAll.Range("F15").Value = 12753
$wsAll.Range("F16").Value = 100
$wsAll.Range("F19").Value = 5246
